# Power tree.xlsx - update DDR3L contribution into the PMIC SW3 rail.
$wb = $excel.ActiveWorkbook

$pmic = $wb.Worksheets.Item("PMIC")
$cpu  = $wb.Worksheets.Item("CPU")
$eth  = $wb.Worksheets.Item("Ethernet PHY")
$ddr  = $wb.Worksheets.Item("DDR3L")

# Row 5 = SW3 rail: its "Sinks" label now also lists the DDR3L sink, and its
# output current (O5) must also account for the DDR3L per-rail current (D8).
$pmic.Range("K5").Value = "NVCC_DRAM, VINREFDDR, DDR3L"
$pmic.Range("O5").Formula = "=CPU!D24+SUMIF(E2:E14, N5, F2:F14)+DDR3L!D8"

# Row 14 = VINREFDDR sink: its output current (O14) must also include the
# DDR3L reference-voltage current (D9).
$pmic.Range("O14").Formula = "=CPU!D25+DDR3L!D9"

# Restore the view state: PMIC becomes the active sheet/tab, scrolled over to
# column J with R27 selected; CPU is no longer the tab-selected sheet; the
# Ethernet PHY sheet's last selection moves to F8.
$eth.Range("F8").Select() | Out-Null
$cpu.Range("B21:B30").Select() | Out-Null
$cpu.Range("B21").Activate() | Out-Null

$pmic.Activate()
$pmic.Range("J1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = $pmic.Range("J1").Column
$pmic.Range("R27").Select() | Out-Null
